$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Araba/Álava', 0.285735240394602, 1),
    @('Albacete', 0.2149145517133375, 1),
    @('Alacant/Alicante', 0.2286231344379116, 1),
    @('Almería', 0.2109161930577687, 1),
    @('Ávila', 0.3005370139837352, 1),
    @('Badajoz', 0.2788425237537643, 1),
    @('Illes Balears', 0.2067202647003389, 1),
    @('Barcelona', 0.2483923808593487, 1),
    @('Burgos', 0.2979317486766824, 1),
    @('Cáceres', 0.2815207392571256, 1),
    @('Cádiz', 0.1597664163184027, 2),
    @('Castelló/Castellón', 0.2761124481894459, 1),
    @('Ciudad Real', 0.2704151871992366, 1),
    @('Córdoba', 0.2856268868679157, 1),
    @('A Coruña', 0.2929625249432078, 1),
    @('Cuenca', 0.3180550172583693, 1),
    @('Girona', 0.2064482205816434, 1),
    @('Granada', 0.2438251478407663, 1),
    @('Guadalajara', 0.298164325351311, 1),
    @('Gipuzkoa', 0.2365932479304993, 1),
    @('Huelva', 0.2929940828484041, 1),
    @('Huesca', 0.2634694327368446, 1),
    @('Jaén', 0.2859912102692259, 1),
    @('León', 0.2928891127375993, 1),
    @('Lleida', 0.242312466940871, 1),
    @('La Rioja', 0.2866321409486296, 1),
    @('Lugo', 0.3095875227688303, 1),
    @('Madrid', 0.2190426832982557, 1),
    @('Málaga', 0.1988923728065914, 1),
    @('Murcia', 0.2433764205369494, 1),
    @('Navarra', 0.165981037437761, 1),
    @('Ourense', 0.2981650117895843, 1),
    @('Asturias', 0.2860200038118403, 1),
    @('Palencia', 0.3105363471715809, 1),
    @('Pontevedra', 0.2957197434628145, 1),
    @('Salamanca', 0.293086470332834, 1),
    @('Cantabria', 0.298987722061843, 1),
    @('Segovia', 0.3018976614970462, 1),
    @('Sevilla', 0.2238729181390343, 1),
    @('Soria', 0.3264392719739629, 1),
    @('Tarragona', 0.230751469302527, 1),
    @('Teruel', 0.2178704947316009, 1),
    @('Toledo', 0.2597829439487325, 1),
    @('València/Valencia', 0.1983907527942733, 1),
    @('Valladolid', 0.2833821107420834, 1),
    @('Bizkaia', 0.293074638430981, 1),
    @('Zamora', 0.292904314164334, 1),
    @('Zaragoza', 0.29216951665411, 1),
    @('Ceuta', 0.0794999015389227, 2),
    @('Melilla', 0.103159321735515, 2)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

# Remove now-obsolete trailing rows 52 and 53 (table shrank from 53 to 51 rows)
$ws.Range("A52:C53").Delete() | Out-Null
